$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D.
# This shifts the existing "Tipo" header (D1->E1) and "single" value (D2->E2)
# one column to the right, matching the diff.
$ws.Range("D1:D2").EntireColumn.Insert()

# Fill in the new "MAE" header in D1, using the same formatting as the
# other header cells (bold, bordered, centered).
$ws.Range("D1").Value = "MAE"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new MAE numeric value in D2.
$ws.Range("D2").Value = 0.3777406617731509

$ws.Application.CutCopyMode = $false
